$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2-18 changes from 45170 (2023-09-01)
# to 45174 (2023-09-05).
$ws.Range("C2:C18").Value = 45174
